$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows to append (sheet rows 189-191, Indice 188-190),
# mirroring the existing row 188 formatting (style on col A + col E).
$rows = @(
    @{
        Row = 189; Indice = 188; Data = 45227.04166666666
        Home = "Huila"; HomeG = 2; Away = "Santa Fe"; AwayG = 2
        J = 2.65; K = "24/10/2023 03:12"; L = 3.55; M = "28/10/2023 00:59"
        N = 3.11; O = "24/10/2023 03:12"; P = 3.4; Q = "28/10/2023 00:59"
        R = 2.9; S = "24/10/2023 03:12"; T = 2.18; U = "28/10/2023 00:59"
        V = "https://www.betexplorer.com/football/colombia/primera-a/huila-santa-fe/fTMI5JAm/"
    },
    @{
        Row = 190; Indice = 189; Data = 45227.13194444445
        Home = "Deportes Tolima"; HomeG = 2; Away = "U. Magdalena"; AwayG = 1
        J = 1.5; K = "23/10/2023 15:12"; L = 1.48; M = "28/10/2023 03:03"
        N = 4.13; O = "23/10/2023 15:12"; P = 4.43; Q = "28/10/2023 03:06"
        R = 7.16; S = "23/10/2023 15:12"; T = 7.31; U = "28/10/2023 03:06"
        V = "https://www.betexplorer.com/football/colombia/primera-a/deportes-tolima-union-magdalena/xCEZ1yAC/"
    },
    @{
        Row = 191; Indice = 190; Data = 45230.08333333334
        Home = "America De Cali"; HomeG = 1; Away = "Millonarios"; AwayG = 0
        J = 2.2; K = "27/10/2023 03:43"; L = 1.71; M = "31/10/2023 01:55"
        N = 3.17; O = "27/10/2023 03:43"; P = 3.66; Q = "31/10/2023 01:59"
        R = 3.66; S = "27/10/2023 03:43"; T = 5.58; U = "31/10/2023 01:59"
        V = "https://www.betexplorer.com/football/colombia/primera-a/america-de-cali-millonarios/4QI6sASS/"
    }
)

$lastRow = 188

foreach ($r in $rows) {
    $row = $r.Row

    # Carry over the formatting used by the previous data row (bold/border
    # style on Indice, date numFmt on data_partida) before writing values.
    $ws.Cells.Item($lastRow, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($lastRow, 5).Copy($ws.Cells.Item($row, 5))

    $ws.Cells.Item($row, 1).Value = $r.Indice
    $ws.Cells.Item($row, 2).Value = "colombia"
    $ws.Cells.Item($row, 3).Value = "primera-a"

    # "2023" parses as a number, so force text formatting to write it as a
    # string, then drop back to the default style (matches neighbouring
    # text cells, which carry no explicit style).
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = "2023"
    $dCell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $r.Data
    $ws.Cells.Item($row, 6).Value = $r.Home
    $ws.Cells.Item($row, 7).Value = $r.HomeG
    $ws.Cells.Item($row, 8).Value = $r.Away
    $ws.Cells.Item($row, 9).Value = $r.AwayG
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V

    $lastRow = $row
}
